$wb = $excel.ActiveWorkbook

function Set-RangeData($ws, $startRow, $endRow, $rowsData) {
    $n = $rowsData.Count
    $cols = $rowsData[0].Count
    $arr = New-Object "object[,]" $n,$cols
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt $cols; $j++) {
            $arr[$i,$j] = $rowsData[$i][$j]
        }
    }
    $endCol = [char](65 + $cols - 1)
    $rangeAddr = "A" + $startRow + ":" + $endCol + $endRow
    $ws.Range($rangeAddr).Value = $arr
}

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets("LP1912")
$ws.Range("A2").Value = "Última actualización: 15:51:48"
$ws.Range("A3").Value = "Total filas: 310"
$rowsData = @(
    @('03:55:13','04:01','81_EL PELIGRO',6,'LP1912'),
    @('03:55:13','04:34','215_ALUAR',39,'LP1912'),
    @('04:18:52','04:45','215A_EL PATO',27,'LP1912'),
    @('03:55:13','04:53','11_ETCHEVERRY',58,'LP1912'),
    @('04:54:17','04:54','11_ETCHEVERRY',0,'LP1912'),
    @('03:55:13','05:16','17_ROMERO',81,'LP1912'),
    @('04:18:52','05:21','23_HERNANDEZ',63,'LP1912'),
    @('03:55:13','05:22','23_HERNANDEZ',87,'LP1912'),
    @('04:18:52','05:34','215B_EL PATO',76,'LP1912'),
    @('03:55:13','05:35','215B_EL PATO',100,'LP1912'),
    @('03:55:13','05:46','15_ABASTO',111,'LP1912'),
    @('04:18:52','05:53','10_OLMOS',95,'LP1912'),
    @('04:40:48','05:54','10_OLMOS',74,'LP1912'),
    @('04:40:48','06:04','16_SANTA ANA',84,'LP1912'),
    @('04:18:52','06:05','16_SANTA ANA',107,'LP1912'),
    @('04:18:52','06:11','215A_EL PATO',113,'LP1912'),
    @('04:18:52','06:13','225_HARAS DEL SUR',115,'LP1912'),
    @('04:40:48','06:14','225_HARAS DEL SUR',94,'LP1912'),
    @('04:40:48','06:21','26_HERNANDEZ',101,'LP1912'),
    @('04:40:48','06:27','23_HERNANDEZ',107,'LP1912'),
    @('04:40:48','06:29','86_EST CHICA-ESC AGRARIA',109,'LP1912'),
    @('05:49:10','06:30','86_EST CHICA-ESC AGRARIA',41,'LP1912'),
    @('04:40:48','06:31','16_SANTA ANA',111,'LP1912'),
    @('04:54:17','06:44','225_C ROCA-H SUR',110,'LP1912'),
    @('04:54:17','06:46','215C_EL PATO',112,'LP1912'),
    @('05:49:10','06:47','215C_EL PATO',58,'LP1912'),
    @('05:19:24','06:59','14_ABASTO',100,'LP1912'),
    @('05:49:10','07:00','14_ABASTO',71,'LP1912'),
    @('06:14:19','07:01','16_SANTA ANA',47,'LP1912'),
    @('05:49:10','07:05','23_HERNANDEZ',76,'LP1912'),
    @('05:19:24','07:05','15_ABASTO',106,'LP1912'),
    @('05:19:24','07:07','225_GOMEZ',108,'LP1912'),
    @('05:19:24','07:11','215A_EL PATO',112,'LP1912'),
    @('05:19:24','07:15','11_ETCHEVERRY',116,'LP1912'),
    @('05:49:10','07:16','11_ETCHEVERRY',87,'LP1912'),
    @('06:43:12','07:17','16_SANTA ANA',34,'LP1912'),
    @('07:19:37','07:19','10_OLMOS',0,'LP1912'),
    @('05:49:10','07:21','26_HERNANDEZ',92,'LP1912'),
    @('06:14:19','07:23','10_OLMOS',69,'LP1912'),
    @('06:14:19','07:31','11_ETCHEVERRY',77,'LP1912'),
    @('07:19:37','07:31','16_SANTA ANA',12,'LP1912'),
    @('05:49:10','07:32','16_SANTA ANA',103,'LP1912'),
    @('05:49:10','07:32','84_COLONIA URQUIZA-ESC 49',103,'LP1912'),
    @('05:49:10','07:32','11_ETCHEVERRY',103,'LP1912'),
    @('07:19:37','07:34','23_HERNANDEZ',15,'LP1912'),
    @('07:19:37','07:36','27_EL RETIRO',17,'LP1912'),
    @('05:49:10','07:37','27_EL RETIRO',108,'LP1912'),
    @('05:49:10','07:39','10_OLMOS',110,'LP1912'),
    @('07:45:49','07:46','16_SANTA ANA',1,'LP1912'),
    @('06:14:19','07:47','14_ABASTO',93,'LP1912'),
    @('05:49:10','07:48','14_ABASTO',119,'LP1912'),
    @('06:14:19','07:51','215D_EL PATO',97,'LP1912'),
    @('06:43:12','07:52','215D_EL PATO',69,'LP1912'),
    @('07:45:49','07:55','10_OLMOS',10,'LP1912'),
    @('07:19:37','07:59','23_HERNANDEZ',40,'LP1912'),
    @('06:14:19','08:00','23_HERNANDEZ',106,'LP1912'),
    @('07:19:37','08:03','11_ETCHEVERRY',44,'LP1912'),
    @('06:43:12','08:03','23_HERNANDEZ',80,'LP1912'),
    @('06:57:11','08:05','23_HERNANDEZ',68,'LP1912'),
    @('07:19:37','08:10','16_SANTA ANA',51,'LP1912'),
    @('07:58:19','08:11','16_SANTA ANA',13,'LP1912'),
    @('06:14:19','08:12','15_ABASTO',118,'LP1912'),
    @('07:45:49','08:13','10_OLMOS',28,'LP1912'),
    @('06:43:12','08:21','26_HERNANDEZ',98,'LP1912'),
    @('07:19:37','08:22','16_P MOR-SANTA ANA',63,'LP1912'),
    @('06:43:12','08:23','215B_EL PATO',100,'LP1912'),
    @('06:43:12','08:23','16_P MOR-SANTA ANA',100,'LP1912'),
    @('06:43:12','08:27','84_COLONIA URQUIZA-ESC 49',104,'LP1912'),
    @('07:45:49','08:31','23_HERNANDEZ',46,'LP1912'),
    @('07:58:19','08:33','10_OLMOS',35,'LP1912'),
    @('07:58:19','08:35','23_HERNANDEZ',37,'LP1912'),
    @('08:36:20','08:38','23_HERNANDEZ',2,'LP1912'),
    @('06:57:11','08:42','81_EL PELIGRO',105,'LP1912'),
    @('07:58:19','08:42','16_SANTA ANA',44,'LP1912'),
    @('08:19:33','08:43','16_SANTA ANA',24,'LP1912'),
    @('07:19:37','08:43','14_ABASTO',84,'LP1912'),
    @('08:36:20','08:46','16_SANTA ANA',10,'LP1912'),
    @('08:48:09','08:49','16_SANTA ANA',1,'LP1912'),
    @('08:19:33','08:53','10_OLMOS',34,'LP1912'),
    @('06:57:11','08:54','17_ROMERO',117,'LP1912'),
    @('08:55:19','08:56','16_SANTA ANA',1,'LP1912'),
    @('08:55:19','08:56','10_OLMOS',1,'LP1912'),
    @('07:19:37','09:01','215A_EL PATO',102,'LP1912'),
    @('08:19:33','09:02','23_HERNANDEZ',43,'LP1912'),
    @('08:48:09','09:02','215A_EL PATO',14,'LP1912'),
    @('07:45:49','09:03','11_ETCHEVERRY',78,'LP1912'),
    @('08:48:09','09:04','11_ETCHEVERRY',16,'LP1912'),
    @('08:36:20','09:05','23_HERNANDEZ',29,'LP1912'),
    @('08:55:19','09:06','23_HERNANDEZ',11,'LP1912'),
    @('07:19:37','09:10','16_P MOR-SANTA ANA',111,'LP1912'),
    @('08:36:20','09:11','16_SANTA ANA',35,'LP1912'),
    @('08:48:09','09:11','16_P MOR-SANTA ANA',23,'LP1912'),
    @('08:48:09','09:13','16_SANTA ANA',25,'LP1912'),
    @('08:36:20','09:13','10_OLMOS',37,'LP1912'),
    @('08:55:19','09:14','16_SANTA ANA',19,'LP1912'),
    @('07:19:37','09:16','27_EL RETIRO',117,'LP1912'),
    @('07:58:19','09:17','27_EL RETIRO',79,'LP1912'),
    @('07:45:49','09:21','26_HERNANDEZ',96,'LP1912'),
    @('07:45:49','09:22','17_ROMERO',97,'LP1912'),
    @('07:58:19','09:23','17_ROMERO',85,'LP1912'),
    @('07:45:49','09:23','11_ETCHEVERRY',98,'LP1912'),
    @('08:48:09','09:24','11_ETCHEVERRY',36,'LP1912'),
    @('08:19:33','09:25','16_SANTA ANA',66,'LP1912'),
    @('07:45:49','09:32','15_ABASTO',107,'LP1912'),
    @('07:45:49','09:33','10_OLMOS',108,'LP1912'),
    @('09:25:56','09:34','23_HERNANDEZ',9,'LP1912'),
    @('08:48:09','09:35','23_HERNANDEZ',47,'LP1912'),
    @('07:45:49','09:41','215C_EL PATO',116,'LP1912'),
    @('07:58:19','09:42','215C_EL PATO',104,'LP1912'),
    @('07:58:19','09:43','14_ABASTO',105,'LP1912'),
    @('08:48:09','09:44','14_ABASTO',56,'LP1912'),
    @('09:25:56','09:52','10_OLMOS',27,'LP1912'),
    @('07:58:19','09:52','15_ABASTO',114,'LP1912'),
    @('09:25:56','09:58','16_SANTA ANA',33,'LP1912'),
    @('09:25:56','10:03','11_ETCHEVERRY',38,'LP1912'),
    @('09:25:56','10:04','23_HERNANDEZ',39,'LP1912'),
    @('08:19:33','10:10','16_P MOR-SANTA ANA',111,'LP1912'),
    @('08:48:09','10:11','16_P MOR-SANTA ANA',83,'LP1912'),
    @('09:25:56','10:12','10_OLMOS',47,'LP1912'),
    @('08:19:33','10:12','15_ABASTO',113,'LP1912'),
    @('09:25:56','10:20','26_HERNANDEZ',55,'LP1912'),
    @('08:36:20','10:21','26_HERNANDEZ',105,'LP1912'),
    @('08:36:20','10:22','17_ROMERO',106,'LP1912'),
    @('09:25:56','10:22','16_SANTA ANA',57,'LP1912'),
    @('09:25:56','10:23','11_ETCHEVERRY',58,'LP1912'),
    @('08:36:20','10:26','215A_EL PATO',110,'LP1912'),
    @('08:48:09','10:27','215A_EL PATO',99,'LP1912'),
    @('10:12:35','10:33','10_OLMOS',21,'LP1912'),
    @('10:12:35','10:34','16_SANTA ANA',22,'LP1912'),
    @('10:12:35','10:35','23_HERNANDEZ',23,'LP1912'),
    @('09:25:56','10:41','17_ROMERO',76,'LP1912'),
    @('08:48:09','10:42','17_ROMERO',114,'LP1912'),
    @('08:55:19','10:43','14_ABASTO',108,'LP1912'),
    @('08:48:09','10:44','14_ABASTO',116,'LP1912'),
    @('10:12:35','10:52','15_ABASTO',40,'LP1912'),
    @('10:52:48','10:52','10_OLMOS',0,'LP1912'),
    @('10:12:35','10:53','10_OLMOS',41,'LP1912'),
    @('10:52:48','10:53','16_SANTA ANA',1,'LP1912'),
    @('09:25:56','10:53','27_EL RETIRO',88,'LP1912'),
    @('10:52:48','10:54','15_ABASTO',2,'LP1912'),
    @('10:12:35','10:56','27_EL RETIRO',44,'LP1912'),
    @('10:52:48','10:56','16_SANTA ANA',4,'LP1912'),
    @('09:25:56','11:01','215C_EL PATO',96,'LP1912'),
    @('10:12:35','11:02','215C_EL PATO',50,'LP1912'),
    @('10:12:35','11:03','11_ETCHEVERRY',51,'LP1912'),
    @('10:12:35','11:04','23_HERNANDEZ',52,'LP1912'),
    @('10:52:48','11:06','16_P MOR-167 Y 521',14,'LP1912'),
    @('09:25:56','11:10','16_P MOR-167 Y 521',105,'LP1912'),
    @('10:52:48','11:11','10_OLMOS',19,'LP1912'),
    @('10:12:35','11:12','15_ABASTO',60,'LP1912'),
    @('11:17:08','11:17','86_EST CHICA-ESC AGRARIA',0,'LP1912'),
    @('11:17:08','11:17','16_SANTA ANA',0,'LP1912'),
    @('11:17:08','11:18','15_ABASTO',1,'LP1912'),
    @('09:25:56','11:19','86_EST CHICA-ESC AGRARIA',114,'LP1912'),
    @('09:25:56','11:20','26_HERNANDEZ',115,'LP1912'),
    @('10:12:35','11:21','26_HERNANDEZ',69,'LP1912'),
    @('10:52:48','11:26','225_C ROCA-H SUR',34,'LP1912'),
    @('10:12:35','11:27','225_C ROCA-H SUR',75,'LP1912'),
    @('10:12:35','11:32','81_EL PELIGRO',80,'LP1912'),
    @('10:52:48','11:34','23_HERNANDEZ',42,'LP1912'),
    @('11:17:08','11:35','23_HERNANDEZ',18,'LP1912'),
    @('10:52:48','11:41','17_ROMERO',49,'LP1912'),
    @('10:12:35','11:42','17_ROMERO',90,'LP1912'),
    @('10:52:48','11:43','10_OLMOS',51,'LP1912'),
    @('11:45:01','11:45','16_SANTA ANA',0,'LP1912'),
    @('11:45:01','11:46','16_SANTA ANA',1,'LP1912'),
    @('11:45:01','11:47','15_ABASTO',2,'LP1912'),
    @('10:12:35','11:51','215B_EL PATO',99,'LP1912'),
    @('10:12:35','11:56','10_OLMOS',104,'LP1912'),
    @('10:52:48','11:58','225_GOMEZ',66,'LP1912'),
    @('10:12:35','11:58','16_P MOR-167 Y 521',106,'LP1912'),
    @('10:12:35','11:59','225_GOMEZ',107,'LP1912'),
    @('11:59:06','11:59','16_SANTA ANA',0,'LP1912'),
    @('11:59:06','12:00','16_SANTA ANA',1,'LP1912'),
    @('11:59:06','12:01','15_ABASTO',2,'LP1912'),
    @('10:12:35','12:04','84_COLONIA URQUIZA-ESC 49',112,'LP1912'),
    @('11:17:08','12:05','23_HERNANDEZ',48,'LP1912'),
    @('10:12:35','12:06','14_ABASTO',114,'LP1912'),
    @('11:45:01','12:06','84_COLONIA URQUIZA-ESC 49',21,'LP1912'),
    @('10:12:35','12:06','16_P MOR-SANTA ANA',114,'LP1912'),
    @('11:17:08','12:07','14_ABASTO',50,'LP1912'),
    @('11:59:06','12:07','84_COLONIA URQUIZA-ESC 49',8,'LP1912'),
    @('11:17:08','12:07','16_P MOR-SANTA ANA',50,'LP1912'),
    @('11:45:01','12:08','23_HERNANDEZ',23,'LP1912'),
    @('11:17:08','12:09','10_OLMOS',52,'LP1912'),
    @('10:52:48','12:09','84_COLONIA URQUIZA-ESC 49',77,'LP1912'),
    @('11:17:08','12:11','84_COLONIA URQUIZA-ESC 49',54,'LP1912'),
    @('11:45:01','12:13','10_OLMOS',28,'LP1912'),
    @('11:59:06','12:13','15_ABASTO',14,'LP1912'),
    @('10:52:48','12:18','14_ABASTO',86,'LP1912'),
    @('10:52:48','12:20','215A_EL PATO',88,'LP1912'),
    @('11:17:08','12:21','215A_EL PATO',64,'LP1912'),
    @('11:17:08','12:21','14_ABASTO',64,'LP1912'),
    @('10:52:48','12:21','26_HERNANDEZ',89,'LP1912'),
    @('10:52:48','12:24','17_ROMERO',92,'LP1912'),
    @('12:27:08','12:27','16_SANTA ANA',0,'LP1912'),
    @('11:45:01','12:35','23_HERNANDEZ',50,'LP1912'),
    @('11:45:01','12:35','11_ETCHEVERRY',50,'LP1912'),
    @('10:52:48','12:36','27_EL RETIRO',104,'LP1912'),
    @('11:17:08','12:37','27_EL RETIRO',80,'LP1912'),
    @('10:52:48','12:38','17_179 Y 38',106,'LP1912'),
    @('11:17:08','12:41','10_OLMOS',84,'LP1912'),
    @('12:27:08','12:47','16_SANTA ANA',20,'LP1912'),
    @('10:52:48','12:48','11_ETCHEVERRY',116,'LP1912'),
    @('11:17:08','12:48','17_ROMERO',91,'LP1912'),
    @('11:17:08','12:49','11_ETCHEVERRY',92,'LP1912'),
    @('12:27:08','12:50','15_ABASTO',23,'LP1912'),
    @('12:54:06','12:55','16_SANTA ANA',1,'LP1912'),
    @('12:27:08','12:55','10_OLMOS',28,'LP1912'),
    @('11:45:01','13:02','15_ABASTO',77,'LP1912'),
    @('12:27:08','13:03','14_ABASTO',36,'LP1912'),
    @('12:54:06','13:05','23_HERNANDEZ',11,'LP1912'),
    @('11:17:08','13:07','16_P MOR-SANTA ANA',110,'LP1912'),
    @('12:27:08','13:08','10_OLMOS',41,'LP1912'),
    @('11:17:08','13:14','215D_EL PATO',117,'LP1912'),
    @('12:54:06','13:15','11_ETCHEVERRY',21,'LP1912'),
    @('11:17:08','13:16','11_ETCHEVERRY',119,'LP1912'),
    @('11:45:01','13:17','17_ROMERO',92,'LP1912'),
    @('11:45:01','13:20','10_OLMOS',95,'LP1912'),
    @('11:45:01','13:21','26_HERNANDEZ',96,'LP1912'),
    @('13:23:09','13:24','16_SANTA ANA',1,'LP1912'),
    @('12:54:06','13:26','15_ABASTO',32,'LP1912'),
    @('11:45:01','13:27','14_ABASTO',102,'LP1912'),
    @('11:59:06','13:31','17_ROMERO',92,'LP1912'),
    @('12:27:08','13:32','10_OLMOS',65,'LP1912'),
    @('12:54:06','13:35','23_HERNANDEZ',41,'LP1912'),
    @('13:23:09','13:36','15_ABASTO',13,'LP1912'),
    @('12:27:08','13:37','23_HERNANDEZ',70,'LP1912'),
    @('12:27:08','13:46','17_ROMERO',79,'LP1912'),
    @('12:54:06','13:47','16_SANTA ANA',53,'LP1912'),
    @('12:54:06','13:51','11_ETCHEVERRY',57,'LP1912'),
    @('11:59:06','13:51','215A_EL PATO',112,'LP1912'),
    @('11:59:06','13:56','225_GOMEZ',117,'LP1912'),
    @('11:59:06','13:57','16_P MOR-167 Y 521',118,'LP1912'),
    @('14:00:52','14:00','14_ABASTO',0,'LP1912'),
    @('14:00:52','14:00','16_SANTA ANA',0,'LP1912'),
    @('14:00:52','14:01','15_ABASTO',1,'LP1912'),
    @('12:27:08','14:04','17_ROMERO',97,'LP1912'),
    @('14:00:52','14:04','23_HERNANDEZ',4,'LP1912'),
    @('13:23:09','14:05','23_HERNANDEZ',42,'LP1912'),
    @('14:00:52','14:06','16_SANTA ANA',6,'LP1912'),
    @('14:00:52','14:16','27_EL RETIRO',16,'LP1912'),
    @('12:27:08','14:17','27_EL RETIRO',110,'LP1912'),
    @('14:00:52','14:19','215C_EL PATO',19,'LP1912'),
    @('14:00:52','14:20','26_HERNANDEZ',20,'LP1912'),
    @('12:27:08','14:20','215C_EL PATO',113,'LP1912'),
    @('12:54:06','14:21','26_HERNANDEZ',87,'LP1912'),
    @('12:54:06','14:39','14_ABASTO',105,'LP1912'),
    @('14:00:52','14:44','14_ABASTO',44,'LP1912'),
    @('14:44:25','14:45','15_ABASTO',1,'LP1912'),
    @('14:00:52','14:56','16_P MOR-SANTA ANA',56,'LP1912'),
    @('13:23:09','14:57','16_P MOR-SANTA ANA',94,'LP1912'),
    @('13:23:09','14:58','215B_EL PATO',95,'LP1912'),
    @('13:23:09','15:00','81_EL PELIGRO',97,'LP1912'),
    @('15:02:32','15:03','15_ABASTO',1,'LP1912'),
    @('14:44:25','15:04','23_HERNANDEZ',20,'LP1912'),
    @('14:00:52','15:04','10_OLMOS',64,'LP1912'),
    @('13:23:09','15:05','10_OLMOS',102,'LP1912'),
    @('14:44:25','15:06','16_SANTA ANA',22,'LP1912'),
    @('14:00:52','15:10','17_ROMERO',70,'LP1912'),
    @('14:00:52','15:13','11_ETCHEVERRY',73,'LP1912'),
    @('13:23:09','15:14','11_ETCHEVERRY',111,'LP1912'),
    @('14:44:25','15:16','16_SANTA ANA',32,'LP1912'),
    @('14:00:52','15:20','15_ABASTO',80,'LP1912'),
    @('13:23:09','15:21','26_HERNANDEZ',118,'LP1912'),
    @('14:00:52','15:25','26_HERNANDEZ',85,'LP1912'),
    @('15:02:32','15:26','16_SANTA ANA',24,'LP1912'),
    @('14:00:52','15:32','84_COLONIA URQUIZA-ESC 49',92,'LP1912'),
    @('15:02:32','15:34','23_HERNANDEZ',32,'LP1912'),
    @('14:00:52','15:35','23_HERNANDEZ',95,'LP1912'),
    @('14:00:52','15:36','10_OLMOS',96,'LP1912'),
    @('14:44:25','15:37','10_OLMOS',53,'LP1912'),
    @('14:00:52','15:38','215A_EL PATO',98,'LP1912'),
    @('14:44:25','15:45','14_ABASTO',61,'LP1912'),
    @('14:00:52','15:46','16_P MOR-167 Y 521',106,'LP1912'),
    @('14:44:25','15:47','23_HERNANDEZ',63,'LP1912'),
    @('14:00:52','15:53','11_ETCHEVERRY',113,'LP1912'),
    @('15:51:48','15:55','17_ROMERO',4,'LP1912'),
    @('14:44:25','15:56','17_ROMERO',72,'LP1912'),
    @('14:00:52','15:56','27_EL RETIRO',116,'LP1912'),
    @('14:44:25','16:01','10_OLMOS',77,'LP1912'),
    @('15:02:32','16:01','27_EL RETIRO',59,'LP1912'),
    @('15:51:48','16:02','16_SANTA ANA',11,'LP1912'),
    @('14:44:25','16:02','27_EL RETIRO',78,'LP1912'),
    @('15:51:48','16:04','23_HERNANDEZ',13,'LP1912'),
    @('15:02:32','16:08','14_ABASTO',66,'LP1912'),
    @('15:51:48','16:10','16_SANTA ANA',19,'LP1912'),
    @('15:51:48','16:14','225_C ROCA-H SUR',23,'LP1912'),
    @('14:44:25','16:15','225_C ROCA-H SUR',91,'LP1912'),
    @('15:51:48','16:19','215C_EL PATO',28,'LP1912'),
    @('14:44:25','16:20','215C_EL PATO',96,'LP1912'),
    @('14:44:25','16:21','26_HERNANDEZ',97,'LP1912'),
    @('15:51:48','16:29','10_OLMOS',38,'LP1912'),
    @('15:51:48','16:30','15_ABASTO',39,'LP1912'),
    @('15:51:48','16:34','23_HERNANDEZ',43,'LP1912'),
    @('15:51:48','16:36','11_ETCHEVERRY',45,'LP1912'),
    @('14:44:25','16:42','16_P MOR-SANTA ANA',118,'LP1912'),
    @('15:51:48','16:42','225_GOMEZ',51,'LP1912'),
    @('14:44:25','16:43','225_GOMEZ',119,'LP1912'),
    @('15:51:48','16:48','15_ABASTO',57,'LP1912'),
    @('15:51:48','16:50','14_ABASTO',59,'LP1912'),
    @('15:02:32','16:56','17_179 Y 38',114,'LP1912'),
    @('15:51:48','17:04','215A_EL PATO',73,'LP1912'),
    @('15:51:48','17:21','26_HERNANDEZ',90,'LP1912'),
    @('15:51:48','17:24','84_COLONIA URQUIZA-ESC 49',93,'LP1912'),
    @('15:51:48','17:28','14_ABASTO',97,'LP1912'),
    @('15:51:48','17:36','27_EL RETIRO',105,'LP1912'),
    @('15:51:48','17:38','17_ROMERO',107,'LP1912'),
    @('15:51:48','17:40','215B_EL PATO',109,'LP1912'),
    @('15:51:48','17:50','16_P MOR-167 Y 521',119,'LP1912')
)
Set-RangeData $ws 6 315 $rowsData

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 15:51:48"
$ws.Range("A3").Value = "Total filas: 32"
$rowsData = @(
    @('03:55:13','04:34','215_ALUAR',39,'LP1912'),
    @('04:18:52','04:45','215A_EL PATO',27,'LP1912'),
    @('04:18:52','05:34','215B_EL PATO',76,'LP1912'),
    @('03:55:13','05:35','215B_EL PATO',100,'LP1912'),
    @('04:18:52','06:11','215A_EL PATO',113,'LP1912'),
    @('04:54:17','06:46','215C_EL PATO',112,'LP1912'),
    @('05:49:10','06:47','215C_EL PATO',58,'LP1912'),
    @('05:19:24','07:11','215A_EL PATO',112,'LP1912'),
    @('06:14:19','07:51','215D_EL PATO',97,'LP1912'),
    @('06:43:12','07:52','215D_EL PATO',69,'LP1912'),
    @('06:43:12','08:23','215B_EL PATO',100,'LP1912'),
    @('07:19:37','09:01','215A_EL PATO',102,'LP1912'),
    @('08:48:09','09:02','215A_EL PATO',14,'LP1912'),
    @('07:45:49','09:41','215C_EL PATO',116,'LP1912'),
    @('07:58:19','09:42','215C_EL PATO',104,'LP1912'),
    @('08:36:20','10:26','215A_EL PATO',110,'LP1912'),
    @('08:48:09','10:27','215A_EL PATO',99,'LP1912'),
    @('09:25:56','11:01','215C_EL PATO',96,'LP1912'),
    @('10:12:35','11:02','215C_EL PATO',50,'LP1912'),
    @('10:12:35','11:51','215B_EL PATO',99,'LP1912'),
    @('10:52:48','12:20','215A_EL PATO',88,'LP1912'),
    @('11:17:08','12:21','215A_EL PATO',64,'LP1912'),
    @('11:17:08','13:14','215D_EL PATO',117,'LP1912'),
    @('11:59:06','13:51','215A_EL PATO',112,'LP1912'),
    @('14:00:52','14:19','215C_EL PATO',19,'LP1912'),
    @('12:27:08','14:20','215C_EL PATO',113,'LP1912'),
    @('13:23:09','14:58','215B_EL PATO',95,'LP1912'),
    @('14:00:52','15:38','215A_EL PATO',98,'LP1912'),
    @('15:51:48','16:19','215C_EL PATO',28,'LP1912'),
    @('14:44:25','16:20','215C_EL PATO',96,'LP1912'),
    @('15:51:48','17:04','215A_EL PATO',73,'LP1912'),
    @('15:51:48','17:40','215B_EL PATO',109,'LP1912')
)
Set-RangeData $ws 6 37 $rowsData

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets("6203-6173")
$ws.Range("A2").Value = "Última actualización: 15:51:48"
$ws.Range("A3").Value = "Total filas: 43"
$rowsData = @(
    @('04:18:52','05:43','215A_LA PLATA',85,'L6173'),
    @('03:55:13','05:44','215A_LA PLATA',109,'L6173'),
    @('04:18:52','06:08','215A_LA PLATA',110,'L6173'),
    @('04:54:17','06:09','215A_LA PLATA',75,'L6173'),
    @('04:40:48','06:32','215C_LA PLATA',112,'L6203'),
    @('04:54:17','06:33','215C_LA PLATA',99,'L6203'),
    @('05:19:24','07:00','215B_LP-P MOR-1 Y 57',101,'L6173'),
    @('05:49:10','07:35','215A_LA PLATA',106,'L6173'),
    @('07:19:37','07:38','215A_LA PLATA',19,'L6173'),
    @('06:14:19','08:07','215C_LA PLATA',113,'L6203'),
    @('07:19:37','08:08','215C_LA PLATA',49,'L6203'),
    @('07:45:49','08:09','215C_LA PLATA',24,'L6203'),
    @('06:57:11','08:10','215C_LA PLATA',73,'L6203'),
    @('07:58:19','08:11','215C_LA PLATA',13,'L6203'),
    @('07:19:37','08:35','215A_LA PLATA',76,'L6173'),
    @('06:57:11','08:38','215A_LA PLATA',101,'L6173'),
    @('07:58:19','08:40','215A_LA PLATA',42,'L6173'),
    @('08:19:33','08:43','215A_LA PLATA',24,'L6173'),
    @('08:48:09','08:53','215A_LA PLATA',5,'L6173'),
    @('08:55:19','08:55','215A_LA PLATA',0,'L6173'),
    @('07:19:37','09:08','215D_LA PLATA',109,'L6203'),
    @('07:58:19','09:09','215D_LA PLATA',71,'L6203'),
    @('08:36:20','09:10','215D_LA PLATA',34,'L6203'),
    @('08:48:09','09:13','215D_LA PLATA',25,'L6203'),
    @('09:25:56','10:02','215B_LP-P MOR-40 Y 115',37,'L6173'),
    @('08:19:33','10:03','215B_LP-P MOR-40 Y 115',104,'L6173'),
    @('09:25:56','10:53','215A_LA PLATA',88,'L6173'),
    @('08:55:19','10:54','215A_LA PLATA',119,'L6173'),
    @('09:25:56','11:13','215C_LA PLATA',108,'L6203'),
    @('10:12:35','12:04','215A_LA PLATA',112,'L6173'),
    @('11:17:08','12:54','215C_LA PLATA',97,'L6203'),
    @('12:54:06','12:55','215C_LA PLATA',1,'L6203'),
    @('11:45:01','13:31','215B_LP-P MOR-1 Y 57',106,'L6173'),
    @('12:27:08','14:09','215A_LA PLATA',102,'L6173'),
    @('14:44:25','14:46','215D_LA PLATA',2,'L6203'),
    @('14:00:52','14:52','215D_LA PLATA',52,'L6203'),
    @('12:54:06','14:53','215D_LA PLATA',119,'L6203'),
    @('14:00:52','15:34','215A_LA PLATA',94,'L6173'),
    @('15:51:48','16:13','215C_LA PLATA',22,'L6203'),
    @('14:44:25','16:14','215C_LA PLATA',90,'L6203'),
    @('15:51:48','16:52','215B_LP-P MOR-40 Y 115',61,'L6173'),
    @('15:02:32','16:53','215B_LP-P MOR-40 Y 115',111,'L6173'),
    @('15:51:48','17:17','215A_LA PLATA',86,'L6173')
)
Set-RangeData $ws 6 48 $rowsData

Write-Output "edit applied"